# Apply odds updates to the active worksheet per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 1.91
$ws.Range("T2").Value = 1.99

# Row 3
$ws.Range("AR3").Value = 1.85
$ws.Range("AS3").Value = 2

# Row 4
$ws.Range("N4").Value = 9.5

# Row 6
$ws.Range("G6").Value = 2.7
$ws.Range("I6").Value = 2.5
$ws.Range("J6").Value = 3.4
$ws.Range("L6").Value = 3.2

$ws.Range("AA6").Value = 9.5
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 26
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 29

$ws.Range("AL6").Value = 9
$ws.Range("AM6").Value = 13
$ws.Range("AN6").Value = 10
$ws.Range("AO6").Value = 26
$ws.Range("AP6").Value = 21
